# Automatic update of files.
# Applies the per-row value changes described by the diff: the sighting
# rows (2-16) on the "Artfynd" sheet were re-exported, which re-ordered a
# few rows (pairwise/cyclic swaps of their content) and bumped several
# "Taxonsorteringsordning" (column B) id values by +1. Column A (Id),
# D..H (species info) and Q/R (coordinates) travel together with each
# physical observation; only B is independently incremented for rows
# whose observation stayed in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (unchanged observation, B bumped)
$ws.Range("B2").Value = 79244

# Row 3 (unchanged observation, B bumped)
$ws.Range("B3").Value = 92268

# Row 4 (now holds what used to be row 5's observation)
$ws.Range("A4").Value = 130834381
$ws.Range("B4").Value = 91809
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 424485
$ws.Range("R4").Value = 6711319
$ws.Range("AF4").Value = ""

# Row 5 (now holds what used to be row 4's observation)
$ws.Range("A5").Value = 130834396
$ws.Range("B5").Value = 92268
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 1209
$ws.Range("F5").Value = "Rynkskinn"
$ws.Range("G5").Value = "Hermanssonia centrifuga"
$ws.Range("H5").Value = "(P. Karst.) Zmitr."
$ws.Range("Q5").Value = 424484
$ws.Range("R5").Value = 6711318

# Row 6 (now holds what used to be row 7's observation)
$ws.Range("A6").Value = 130834375
$ws.Range("B6").Value = 91772
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 5447
$ws.Range("F6").Value = "Vedticka"
$ws.Range("G6").Value = "Fuscoporia viticola"
$ws.Range("H6").Value = "(Schwein.) Murrill"
$ws.Range("Q6").Value = 424498
$ws.Range("R6").Value = 6711351
$ws.Range("AF6").Value = ""

# Row 7 (now holds what used to be row 6's observation)
$ws.Range("A7").Value = 130834394
$ws.Range("B7").Value = 92268
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 1209
$ws.Range("F7").Value = "Rynkskinn"
$ws.Range("G7").Value = "Hermanssonia centrifuga"
$ws.Range("H7").Value = "(P. Karst.) Zmitr."
$ws.Range("Q7").Value = 424506
$ws.Range("R7").Value = 6711370

# Row 8 (now holds what used to be row 9's observation)
$ws.Range("A8").Value = 130834377
$ws.Range("B8").Value = 91809
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = "Ullticka"
$ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q8").Value = 424489
$ws.Range("R8").Value = 6711391
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = ""

# Row 9 (now holds what used to be row 8's observation)
$ws.Range("A9").Value = 130834387
$ws.Range("B9").Value = 57884
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("Q9").Value = 424517
$ws.Range("R9").Value = 6711378
$ws.Range("M9").Value = "äldre spår"

# Row 10 (unchanged observation, B bumped)
$ws.Range("B10").Value = 91809

# Row 11 (now holds what used to be row 12's observation)
$ws.Range("A11").Value = 130834378
$ws.Range("B11").Value = 91809
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = "Ullticka"
$ws.Range("G11").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q11").Value = 424513
$ws.Range("R11").Value = 6711372
$ws.Range("AF11").Value = ""

# Row 12 (now holds what used to be row 13's observation)
$ws.Range("A12").Value = 130834395
$ws.Range("B12").Value = 92268
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 1209
$ws.Range("F12").Value = "Rynkskinn"
$ws.Range("G12").Value = "Hermanssonia centrifuga"
$ws.Range("H12").Value = "(P. Karst.) Zmitr."
$ws.Range("Q12").Value = 424495
$ws.Range("R12").Value = 6711339

# Row 13 (now holds what used to be row 11's observation)
$ws.Range("A13").Value = 130834392
$ws.Range("B13").Value = 92268
$ws.Range("Q13").Value = 424488
$ws.Range("R13").Value = 6711452

# Row 14 (unchanged observation, B bumped)
$ws.Range("B14").Value = 91809

# Row 15 (unchanged observation, B bumped)
$ws.Range("B15").Value = 91809

# Row 16 (unchanged observation, B bumped)
$ws.Range("B16").Value = 91772
